$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> Computer Company value for column O
$companyByRow = @{
    2 = "Sony"
    3 = "Sony"
    4 = "Sony"
    5 = "Sony"
    6 = "Sony"
    7 = "Sony"
    8 = "Sony"
    9 = "Sony"
    10 = "Sony"
    11 = "Sony"
    12 = "Sony"
    13 = "Sony"
    14 = "Sony"
    15 = "Apple Inc."
    16 = "Apple Inc."
    17 = "Apple Inc."
    18 = "Apple Inc."
    19 = "Apple Inc."
    20 = "Apple Inc."
    21 = "System76"
    22 = "System76"
    23 = "System76"
    24 = "System76"
    25 = "System76"
    26 = "System76"
    27 = "Sony"
    28 = "Sony"
    29 = "Sony"
    30 = "Sony"
    31 = "Sony"
    32 = "Sony"
    33 = "Apple Inc."
    34 = "Apple Inc."
    35 = "Sony"
    36 = "Sony"
    37 = "Sony"
    38 = "Sony"
    39 = "Sony"
    40 = "Sony"
    41 = "Sony"
    42 = "Sony"
    43 = "Sony"
    44 = "Sony"
    45 = "Sony"
    46 = "Sony"
    47 = "Apple Inc."
    48 = "Apple Inc."
    49 = "System76"
    50 = "System76"
    51 = "Sony"
    52 = "IBM"
    53 = "IBM"
    54 = "IBM"
    55 = "IBM"
    56 = "Sony"
    57 = "Sony"
    58 = "Sony"
    59 = "Sony"
    60 = "Sony"
    61 = "Sony"
    62 = "Sony"
    63 = "Sony"
    64 = "Sony"
    65 = "Sony"
    66 = "Sony"
    67 = "Sony"
    68 = "Sony"
    69 = "Sony"
    70 = "Sony"
    71 = "Sony"
}

foreach ($row in $companyByRow.Keys) {
    $ws.Range("O" + $row).Value = $companyByRow[$row]
}
